$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add new "common part" problem 13 content (rows 30-34), interleaved with the
#     touch-ups to existing descriptions in the order the author made them ---
$ws.Range("A30").Value = "c0020"
$ws.Range("B30").Value = "두 점을 지나는 직선의 방정식을 구합니다."

$ws.Range("A31").Value = "c0021"

# Existing description reused/updated to use $x$ math markup
$ws.Range("B50").Value = "둘러싸인 부분의 넓이를 구하기 위해 두 식을 연립해서 교점의 `$x`$-성분을 구합니다. "

$ws.Range("B31").Value = "직선의 방정식에 `$x=0`$을 대입하여 `$y`$ 절편을 구합니다. "

$ws.Range("A32").Value = "c0022"
$ws.Range("B32").Value = "두 점의 `$y`$ 좌표가 일치하도록 방정식을 세웁니다."

$ws.Range("A33").Value = "c0023"
$ws.Range("B34").Value = "미정계수가 포함된 함수식에 주어진 `$x`$ 값을 대입해서 미정계수 사이의 관계식을 구합니다."

# Existing description reused/updated to include $f(2)$
$ws.Range("C56").Value = "구간함수의 함숫값; `$f(2)`$"

$ws.Range("C30").Value = "`$\left(a, \log _{2} a\right),\left(b, \log _{2} b\right)`$; `$\left(a, \log _{4} a\right),\left(b, \log _{4} b\right)`$; "

$ws.Range("A34").Value = "c0024"
$ws.Range("B33").Value = "로그방정식을 풀기 위해 밑을 통일하고 로그법칙을 이용해서 변변 정리해줍니다."
$ws.Range("C33").Value = "결과 `$a^{b}=b^{a}`$;"
$ws.Range("C34").Value = "조건 `$f(1)=40`$;"

# --- Update the sheet view: drop the frozen top-left scroll position and move the
#     selection to the newly added C34 cell ---
$ws.Activate()
$ws.Range("C34").Select()
